$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Dosage")

# Add the new PlennyShake version note next to the existing history list (G2:G4)
$ws.Range("G5").Value = "PlennyShake 2020 (3.0) : 416 calories / 100 grammes, sachet 960 grammes"

# Correction d'une erreur de copie dans une cellule (C18/C19 pointaient sur B4/B5 au lieu de B3)
$ws.Range("C18").Formula = "=B18/(B3/100)"
$ws.Range("C19").Formula = "=B19/(B3/100)"

# Restore the selection left by the author after the edit
$ws.Range("G6").Select()
